# update on 20210731 画中人
# Replace surrounding double quotes with single quotes in a handful of
# English story lines (EN column = C) in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C17").Value = "[name=`"???`"]   So that's the 'Rhodes Island' that I've been hearing so much about? Heh... How heroic.`n"
$ws.Range("C35").Value = "[name=`"Nearl`"]  'Boiling dark clouds, tumbling in flames...'`n"
$ws.Range("C39").Value = "[name=`"???`"]  '...As terror strips away their voices, the land falls into silence.'`n"
$ws.Range("C40").Value = "[name=`"???`"]  'The titanic Originium lowers its head...'`n"
$ws.Range("C41").Value = "[name=`"???`"]  '...and falls upon the scorched shadow of death.'`n"
